$wb = $excel.ActiveWorkbook

# --- Sheet1 (Train Results) ---
$ws1 = $wb.Worksheets.Item(1)

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 40
$row[0,2] = 4
$row[0,3] = 0
$row[0,4] = 28
$row[0,5] = 4
$row[0,6] = 4
$row[0,7] = 20
$row[0,8] = 2.31
$row[0,9] = 2.593070030212402
$ws1.Range("A2:J2").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 4
$row[0,1] = 0
$row[0,2] = 8
$row[0,3] = 24
$row[0,4] = 4
$row[0,5] = 16
$row[0,6] = 44
$row[0,7] = 0
$row[0,8] = 2.98
$row[0,9] = 3.10846471786499
$ws1.Range("A3:J3").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 4
$row[0,1] = 16
$row[0,2] = 0
$row[0,3] = 20
$row[0,4] = 4
$row[0,5] = 12
$row[0,6] = 44
$row[0,7] = 0
$row[0,8] = 3.22
$row[0,9] = 3.164594650268555
$ws1.Range("A4:J4").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 24
$row[0,2] = 4
$row[0,3] = 0
$row[0,4] = 8
$row[0,5] = 4
$row[0,6] = 36
$row[0,7] = 24
$row[0,8] = 2.45
$row[0,9] = 2.393349647521973
$ws1.Range("A5:J5").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 4
$row[0,1] = 0
$row[0,2] = 8
$row[0,3] = 4
$row[0,4] = 24
$row[0,5] = 12
$row[0,6] = 44
$row[0,7] = 4
$row[0,8] = 3.17
$row[0,9] = 3.100662231445312
$ws1.Range("A6:J6").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 16
$row[0,2] = 8
$row[0,3] = 4
$row[0,4] = 16
$row[0,5] = 20
$row[0,6] = 36
$row[0,7] = 0
$row[0,8] = 3.07
$row[0,9] = 2.920608997344971
$ws1.Range("A7:J7").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 4
$row[0,1] = 20
$row[0,2] = 4
$row[0,3] = 4
$row[0,4] = 16
$row[0,5] = 0
$row[0,6] = 52
$row[0,7] = 0
$row[0,8] = 3.35
$row[0,9] = 3.10871696472168
$ws1.Range("A8:J8").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 4
$row[0,1] = 0
$row[0,2] = 12
$row[0,3] = 4
$row[0,4] = 4
$row[0,5] = 16
$row[0,6] = 52
$row[0,7] = 8
$row[0,8] = 2.81
$row[0,9] = 2.955140829086304
$ws1.Range("A9:J9").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 4
$row[0,1] = 12
$row[0,2] = 8
$row[0,3] = 0
$row[0,4] = 8
$row[0,5] = 8
$row[0,6] = 56.00000000000001
$row[0,7] = 4
$row[0,8] = 2.88
$row[0,9] = 2.931411504745483
$ws1.Range("A10:J10").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 4
$row[0,1] = 4
$row[0,2] = 4
$row[0,3] = 8
$row[0,4] = 28
$row[0,5] = 16
$row[0,6] = 32
$row[0,7] = 4
$row[0,8] = 3.04
$row[0,9] = 3.164654493331909
$ws1.Range("A11:J11").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 4
$row[0,1] = 0
$row[0,2] = 8
$row[0,3] = 4
$row[0,4] = 24
$row[0,5] = 12
$row[0,6] = 44
$row[0,7] = 4
$row[0,8] = 3.09
$row[0,9] = 3.100662231445312
$ws1.Range("A12:J12").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 4
$row[0,1] = 8
$row[0,2] = 0
$row[0,3] = 12
$row[0,4] = 20
$row[0,5] = 8
$row[0,6] = 48
$row[0,7] = 0
$row[0,8] = 3.11
$row[0,9] = 3.159119844436646
$ws1.Range("A13:J13").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 20
$row[0,2] = 8
$row[0,3] = 4
$row[0,4] = 0
$row[0,5] = 20
$row[0,6] = 44
$row[0,7] = 4
$row[0,8] = 2.82
$row[0,9] = 2.848893165588379
$ws1.Range("A14:J14").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 20
$row[0,2] = 8
$row[0,3] = 4
$row[0,4] = 0
$row[0,5] = 20
$row[0,6] = 44
$row[0,7] = 4
$row[0,8] = 2.74
$row[0,9] = 2.848893165588379
$ws1.Range("A15:J15").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 4
$row[0,1] = 12
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 16
$row[0,5] = 8
$row[0,6] = 52
$row[0,7] = 8
$row[0,8] = 2.62
$row[0,9] = 2.861652135848999
$ws1.Range("A16:J16").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 4
$row[0,1] = 12
$row[0,2] = 8
$row[0,3] = 4
$row[0,4] = 16
$row[0,5] = 8
$row[0,6] = 48
$row[0,7] = 0
$row[0,8] = 3.19
$row[0,9] = 3.103657007217407
$ws1.Range("A17:J17").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 20
$row[0,2] = 4
$row[0,3] = 0
$row[0,4] = 4
$row[0,5] = 4
$row[0,6] = 48
$row[0,7] = 20
$row[0,8] = 2.52
$row[0,9] = 2.408069133758545
$ws1.Range("A18:J18").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 4
$row[0,1] = 8
$row[0,2] = 4
$row[0,3] = 4
$row[0,4] = 16
$row[0,5] = 12
$row[0,6] = 48
$row[0,7] = 4
$row[0,8] = 2.94
$row[0,9] = 3.017695188522339
$ws1.Range("A19:J19").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 4
$row[0,1] = 12
$row[0,2] = 0
$row[0,3] = 12
$row[0,4] = 20
$row[0,5] = 8
$row[0,6] = 44
$row[0,7] = 0
$row[0,8] = 3.29
$row[0,9] = 3.179691553115845
$ws1.Range("A20:J20").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 4
$row[0,1] = 12
$row[0,2] = 0
$row[0,3] = 12
$row[0,4] = 20
$row[0,5] = 8
$row[0,6] = 44
$row[0,7] = 0
$row[0,8] = 3.25
$row[0,9] = 3.179691553115845
$ws1.Range("A21:J21").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 12
$row[0,2] = 8
$row[0,3] = 4
$row[0,4] = 16
$row[0,5] = 20
$row[0,6] = 36
$row[0,7] = 4
$row[0,8] = 3.01
$row[0,9] = 2.897711992263794
$ws1.Range("A22:J22").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 4
$row[0,1] = 12
$row[0,2] = 4
$row[0,3] = 4
$row[0,4] = 20
$row[0,5] = 16
$row[0,6] = 32
$row[0,7] = 8
$row[0,8] = 2.96
$row[0,9] = 3.032543182373047
$ws1.Range("A23:J23").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 4
$row[0,1] = 8
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 4
$row[0,5] = 8
$row[0,6] = 52
$row[0,7] = 24
$row[0,8] = 2.44
$row[0,9] = 2.457048416137695
$ws1.Range("A24:J24").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 4
$row[0,1] = 0
$row[0,2] = 4
$row[0,3] = 0
$row[0,4] = 20
$row[0,5] = 8
$row[0,6] = 52
$row[0,7] = 12
$row[0,8] = 2.92
$row[0,9] = 2.912370204925537
$ws1.Range("A25:J25").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 4
$row[0,1] = 8
$row[0,2] = 12
$row[0,3] = 4
$row[0,4] = 4
$row[0,5] = 24
$row[0,6] = 39.99999999999999
$row[0,7] = 4
$row[0,8] = 2.9
$row[0,9] = 3.012637615203857
$ws1.Range("A26:J26").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 4
$row[0,1] = 0
$row[0,2] = 12
$row[0,3] = 16
$row[0,4] = 4
$row[0,5] = 12
$row[0,6] = 52
$row[0,7] = 0
$row[0,8] = 2.96
$row[0,9] = 3.038510799407959
$ws1.Range("A27:J27").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 16
$row[0,2] = 8
$row[0,3] = 4
$row[0,4] = 16
$row[0,5] = 20
$row[0,6] = 36
$row[0,7] = 0
$row[0,8] = 3.02
$row[0,9] = 2.920608997344971
$ws1.Range("A28:J28").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 4
$row[0,1] = 12
$row[0,2] = 8
$row[0,3] = 0
$row[0,4] = 8
$row[0,5] = 8
$row[0,6] = 56.00000000000001
$row[0,7] = 4
$row[0,8] = 2.79
$row[0,9] = 2.931411504745483
$ws1.Range("A29:J29").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 4
$row[0,1] = 0
$row[0,2] = 4
$row[0,3] = 4
$row[0,4] = 16
$row[0,5] = 12
$row[0,6] = 56.00000000000001
$row[0,7] = 4
$row[0,8] = 2.85
$row[0,9] = 2.994871854782104
$ws1.Range("A30:J30").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 4
$row[0,2] = 8
$row[0,3] = 4
$row[0,4] = 16
$row[0,5] = 20
$row[0,6] = 44
$row[0,7] = 4
$row[0,8] = 2.94
$row[0,9] = 2.909976243972778
$ws1.Range("A31:J31").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 12
$row[0,2] = 4
$row[0,3] = 0
$row[0,4] = 4
$row[0,5] = 16
$row[0,6] = 39.99999999999999
$row[0,7] = 24
$row[0,8] = 2.51
$row[0,9] = 2.498975276947021
$ws1.Range("A32:J32").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 8
$row[0,3] = 4
$row[0,4] = 16
$row[0,5] = 12
$row[0,6] = 52
$row[0,7] = 8
$row[0,8] = 2.99
$row[0,9] = 2.936014890670776
$ws1.Range("A33:J33").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 16
$row[0,2] = 0
$row[0,3] = 4
$row[0,4] = 20
$row[0,5] = 20
$row[0,6] = 28
$row[0,7] = 12
$row[0,8] = 3.38
$row[0,9] = 2.88578987121582
$ws1.Range("A34:J34").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 4
$row[0,1] = 16
$row[0,2] = 8
$row[0,3] = 0
$row[0,4] = 12
$row[0,5] = 0
$row[0,6] = 48
$row[0,7] = 12
$row[0,8] = 2.56
$row[0,9] = 2.756472110748291
$ws1.Range("A35:J35").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 8
$row[0,3] = 4
$row[0,4] = 16
$row[0,5] = 12
$row[0,6] = 52
$row[0,7] = 8
$row[0,8] = 2.82
$row[0,9] = 2.936014890670776
$ws1.Range("A36:J36").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 4
$row[0,1] = 12
$row[0,2] = 4
$row[0,3] = 4
$row[0,4] = 12
$row[0,5] = 16
$row[0,6] = 36
$row[0,7] = 12
$row[0,8] = 2.86
$row[0,9] = 2.901228427886963
$ws1.Range("A37:J37").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 24
$row[0,2] = 8
$row[0,3] = 0
$row[0,4] = 12
$row[0,5] = 16
$row[0,6] = 32
$row[0,7] = 8
$row[0,8] = 2.93
$row[0,9] = 2.866098165512085
$ws1.Range("A38:J38").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 4
$row[0,1] = 32
$row[0,2] = 8
$row[0,3] = 4
$row[0,4] = 4
$row[0,5] = 20
$row[0,6] = 24
$row[0,7] = 4
$row[0,8] = 2.84
$row[0,9] = 2.993147373199463
$ws1.Range("A39:J39").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 24
$row[0,2] = 8
$row[0,3] = 0
$row[0,4] = 12
$row[0,5] = 16
$row[0,6] = 32
$row[0,7] = 8
$row[0,8] = 2.94
$row[0,9] = 2.866098165512085
$ws1.Range("A40:J40").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 4
$row[0,1] = 40
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 12
$row[0,5] = 4
$row[0,6] = 36
$row[0,7] = 4
$row[0,8] = 3.16
$row[0,9] = 3.089262008666992
$ws1.Range("A41:J41").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 16
$row[0,2] = 0
$row[0,3] = 4
$row[0,4] = 20
$row[0,5] = 20
$row[0,6] = 28
$row[0,7] = 12
$row[0,8] = 2.72
$row[0,9] = 2.88578987121582
$ws1.Range("A42:J42").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 4
$row[0,1] = 20
$row[0,2] = 4
$row[0,3] = 4
$row[0,4] = 16
$row[0,5] = 0
$row[0,6] = 52
$row[0,7] = 0
$row[0,8] = 3.21
$row[0,9] = 3.10871696472168
$ws1.Range("A43:J43").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 4
$row[0,1] = 12
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 16
$row[0,5] = 8
$row[0,6] = 52
$row[0,7] = 8
$row[0,8] = 2.56
$row[0,9] = 2.86165189743042
$ws1.Range("A44:J44").Value2 = $row

# --- Sheet2 (Test Results) ---
$ws2 = $wb.Worksheets.Item(2)

$row = New-Object 'object[,]' 1,10
$row[0,0] = 4
$row[0,1] = 0
$row[0,2] = 12
$row[0,3] = 16
$row[0,4] = 4
$row[0,5] = 12
$row[0,6] = 52
$row[0,7] = 0
$row[0,8] = 2.96
$row[0,9] = 3.038510799407959
$ws2.Range("A2:J2").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 4
$row[0,1] = 40
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 12
$row[0,5] = 4
$row[0,6] = 36
$row[0,7] = 4
$row[0,8] = 3.16
$row[0,9] = 3.089262008666992
$ws2.Range("A3:J3").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 4
$row[0,1] = 12
$row[0,2] = 4
$row[0,3] = 4
$row[0,4] = 12
$row[0,5] = 16
$row[0,6] = 36
$row[0,7] = 12
$row[0,8] = 2.86
$row[0,9] = 2.901228427886963
$ws2.Range("A4:J4").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 20
$row[0,2] = 8
$row[0,3] = 4
$row[0,4] = 0
$row[0,5] = 20
$row[0,6] = 44
$row[0,7] = 4
$row[0,8] = 2.82
$row[0,9] = 2.848893165588379
$ws2.Range("A5:J5").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 4
$row[0,1] = 12
$row[0,2] = 8
$row[0,3] = 4
$row[0,4] = 16
$row[0,5] = 8
$row[0,6] = 48
$row[0,7] = 0
$row[0,8] = 3.19
$row[0,9] = 3.103657007217407
$ws2.Range("A6:J6").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 4
$row[0,1] = 0
$row[0,2] = 8
$row[0,3] = 4
$row[0,4] = 24
$row[0,5] = 12
$row[0,6] = 44
$row[0,7] = 4
$row[0,8] = 3.09
$row[0,9] = 3.100662231445312
$ws2.Range("A7:J7").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 4
$row[0,2] = 8
$row[0,3] = 4
$row[0,4] = 16
$row[0,5] = 20
$row[0,6] = 44
$row[0,7] = 4
$row[0,8] = 2.94
$row[0,9] = 2.909976243972778
$ws2.Range("A8:J8").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 4
$row[0,1] = 12
$row[0,2] = 0
$row[0,3] = 12
$row[0,4] = 20
$row[0,5] = 8
$row[0,6] = 44
$row[0,7] = 0
$row[0,8] = 3.25
$row[0,9] = 3.179691553115845
$ws2.Range("A9:J9").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 20
$row[0,2] = 8
$row[0,3] = 4
$row[0,4] = 0
$row[0,5] = 20
$row[0,6] = 44
$row[0,7] = 4
$row[0,8] = 2.74
$row[0,9] = 2.848893165588379
$ws2.Range("A10:J10").Value2 = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 4
$row[0,1] = 8
$row[0,2] = 4
$row[0,3] = 4
$row[0,4] = 16
$row[0,5] = 12
$row[0,6] = 48
$row[0,7] = 4
$row[0,8] = 2.94
$row[0,9] = 3.017695188522339
$ws2.Range("A11:J11").Value2 = $row

Write-Output "Edit complete"